$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data row was recorded ahead of the existing row 81, shifting the
# previous rows 81:200 down to 82:201 (dimension grows from R200 to R201).
$ws.Rows("81:81").Insert()

# Populate the newly inserted row 81 with its data.
$ws.Range("A81").Value = 10
$ws.Range("B81").Value = 'Vega Modelo de Temuco'
$ws.Range("C81").Value = 'La Araucanía'
$ws.Range("D81").Value = 44477
$ws.Range("E81").Value = 9
$ws.Range("F81").Value = 100112009
$ws.Range("G81").Value = 'Acelga'
$ws.Range("H81").Value = 'Sin especificar'
$ws.Range("I81").Value = 'Primera'
$ws.Range("J81").Value = 20
$ws.Range("K81").Value = 8000
$ws.Range("L81").Value = 8000
$ws.Range("M81").Value = 8000
$ws.Range("N81").Value = '$/docena de atados (12 kilos)'
$ws.Range("O81").Value = 'Provincia de Cautín'
$ws.Range("P81").Value = 667
$ws.Range("Q81").Value = 12
$ws.Range("R81").Value = 'Hortaliza'
